$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.220.09"
$ws.Range("E2").Value = "  +4.06%  "
$ws.Range("D3").Value = "3.048.10"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.31"
$ws.Range("E5").Value = "  +2.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.02"
$ws.Range("E6").Value = "  +3.31%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "3.038.13"
$ws.Range("E8").Value = "  +1.90%  "
$ws.Range("E9").Value = "  +1.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.137"
$ws.Range("E10").Value = "  +4.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.42"
$ws.Range("E11").Value = "  +11.91%  "
$ws.Range("E12").Value = "  +1.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000236"
$ws.Range("E13").Value = "  +2.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.60"
$ws.Range("E14").Value = "  +2.58%  "
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("D16").Value = "3.554.48"
$ws.Range("E16").Value = "  +2.36%  "
$ws.Range("E17").Value = "  +3.07%  "
$ws.Range("D18").Value = "3.047.55"
$ws.Range("E18").Value = "  +2.17%  "
$ws.Range("D19").Value = "61.116.14"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "445.51"
$ws.Range("E20").Value = "  +4.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.84"
$ws.Range("E21").Value = "  +2.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.730"
$ws.Range("E22").Value = "  +2.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.23"
$ws.Range("E23").Value = "  +1.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.58"
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.60"
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.27"
$ws.Range("E27").Value = "  +8.39%  "
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("E29").Value = "  +3.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.04"
$ws.Range("E30").Value = "  +4.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.42"
$ws.Range("E31").Value = "  +5.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.34"
$ws.Range("E32").Value = "  +2.52%  "
$ws.Range("E33").Value = "  +6.58%  "
$ws.Range("D34").Value = "0.0₃0804"
$ws.Range("E34").Value = "  +6.32%  "
$ws.Range("E35").Value = "  +3.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.02"
$ws.Range("E36").Value = "  +4.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.15"
$ws.Range("E37").Value = "  +3.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "50.02"
$ws.Range("E38").Value = "  +3.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.93"
$ws.Range("E39").Value = "  +6.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.76"
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "410.67"
$ws.Range("E41").Value = "  +3.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0361"
$ws.Range("E42").Value = "  +4.16%  "
$ws.Range("D43").Value = "2.783.26"
$ws.Range("E43").Value = "  +2.20%  "
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("E45").Value = "  +7.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "37.38"
$ws.Range("E46").Value = "  +16.57%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.08"
$ws.Range("E48").Value = "  +3.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.21"
$ws.Range("E49").Value = "  -1.83%  "
$ws.Range("E50").Value = "  +1.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.96"
$ws.Range("E51").Value = "  +2.86%  "
